$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '63.638.30'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -1.07%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.398.38'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("E4").Value = '  -0.08%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '567.10'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '156.80'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.58%  '

$ws.Range("E7").Value = '  -0.06%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '3.400.21'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.12%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.568'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -6.05%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '7.19'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +0.06%  '

$ws.Range("E11").Value = '  -2.05%  '

$ws.Range("E12").Value = '  -3.22%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '3.984.10'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -0.13%  '

$ws.Range("E14").Value = '  -0.33%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '26.96'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.49%  '

$ws.Range("E16").Value = '  -7.46%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '63.706.79'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.06%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.363.79'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.16%  '

$ws.Range("E19").Value = '  -3.42%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '13.53'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -2.80%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '375.47'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.40%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '7.68'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -3.29%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.06%  '

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '70.85'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.55%  '

$ws.Range("E25").Value = '  -5.67%  '

$ws.Range("E26").Value = '  -1.62%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '9.66'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -2.13%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '0.177'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.86%  '

$ws.Range("E29").Value = '  -0.13%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '6.01'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.72%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.38'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -5.11%  '

$ws.Range("E32").Value = '  -2.07%  '

$ws.Range("E33").Value = '  +0.05%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '22.75'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.61%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '6.93'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -2.83%  '

$ws.Range("E36").Value = '  -2.95%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '160.21'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.19%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.83'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.38%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.821'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +7.43%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '25.98'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -2.22%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.0725'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -4.33%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.788.04'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.39%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '42.58'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.32%  '

$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '4.39'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -5.41%  '

$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '6.34'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -6.84%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '25.48'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.15%  '

$ws.Range("E47").Value = '  -2.88%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.33'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +10.07%  '

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '322.09'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +4.06%  '

$ws.Range("E50").Value = '  -3.51%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '6.30'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -3.57%  '
